# Add search criteria ("Restriction") column to the OSS list template.
# Insert a new column before the existing "Obligation" column (H) and
# label its header, shifting subsequent columns one place to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("H:H").Insert()
$ws.Range("H1").Value = "Restriction"
